$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item("Sheet1")

# Create the new "DTYEdit" worksheet right after Sheet1
$ws = $wb.Worksheets.Add($null, $sheet1)
$ws.Name = "DTYEdit"

# Header row - reuse Sheet1's header formatting (style used by A1:B1)
$sheet1.Range("A1:B1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

# Data row - reuse Sheet1's data-row formatting (style used by A2:B2)
$sheet1.Range("A2:B2").Copy()
$ws.Range("A2:B2").PasteSpecial(-4122)

$ws.Range("A1").Value = "orderType"
$ws.Range("B1").Value = "noOfOrders"
$ws.Range("A2").Value = "showOrder"
$ws.Range("B2").Value = 2

$ws.Rows.Item(1).RowHeight = 15
$ws.Rows.Item(2).RowHeight = 15

$ws.Columns.Item(1).ColumnWidth = 9.33
$ws.Columns.Item(2).ColumnWidth = 10.67

# Leave Sheet1's own selection on A1:C2 (no active-sheet highlight there anymore)
$sheet1.Range("A1:C2").Select()

# DTYEdit is the newly active/selected sheet, with G9 selected
$ws.Activate()
$ws.Range("G9").Select()
